# Reorganizing test default models and unit models, refactoring test methods.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheets _set_SUP_TECH -> _set_TECH_SUP, _set_DEM_TECH -> _set_TECH_DEM
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("_set_SUP_TECH").Name = "_set_TECH_SUP"
$wb.Worksheets.Item("_set_DEM_TECH").Name = "_set_TECH_DEM"

# ---------------------------------------------------------------------------
# 2) _set_TECH_DEM: "Final demand" -> "households"
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("_set_TECH_DEM").Range("A2").Value = "households"

# ---------------------------------------------------------------------------
# 3) _set_FLOWS: rename header f_Name_agg -> f_agg_Name, drop the f_Unit column
# ---------------------------------------------------------------------------
$wsFlows = $wb.Worksheets.Item("_set_FLOWS")
$wsFlows.Range("B1").Value = "f_agg_Name"
$wsFlows.Range("C1:C4").Delete()

# ---------------------------------------------------------------------------
# 3b) _set_FLOWS_AGG: header fa_Name -> f_agg_Name
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("_set_FLOWS_AGG").Range("A1").Value = "f_agg_Name"

# ---------------------------------------------------------------------------
# 4) Insert three new sheets after _set_FLOWS_AGG (before _set_YEARS):
#    _set_COSTS, _set_EMISSIONS, _set_LOADFACTORS
# ---------------------------------------------------------------------------
$wsFlowsAgg = $wb.Worksheets.Item("_set_FLOWS_AGG")
$wsCosts = $wb.Worksheets.Add($null, $wsFlowsAgg)
$wsCosts.Name = "_set_COSTS"

$wsEmissions = $wb.Worksheets.Add($null, $wb.Worksheets.Item("_set_COSTS"))
$wsEmissions.Name = "_set_EMISSIONS"

$wsLoadFactors = $wb.Worksheets.Add($null, $wb.Worksheets.Item("_set_EMISSIONS"))
$wsLoadFactors.Name = "_set_LOADFACTORS"

# _set_COSTS: two-column header c_Name / c_Category, two data rows
$wb.Worksheets.Item("_set_FLOWS").Range("A1:B1").Copy($wb.Worksheets.Item("_set_COSTS").Range("A1:B1"))
$wsCosts = $wb.Worksheets.Item("_set_COSTS")
$wsCosts.Range("A1").Value = "c_Name"
$wsCosts.Range("B1").Value = "c_Category"
$wsCosts.Range("A2").Value = "Capital costs"
$wsCosts.Range("B2").Value = "Capital costs"
$wsCosts.Range("A3").Value = "Operational costs"
$wsCosts.Range("B3").Value = "Operational costs"

# _set_EMISSIONS: single-column header e_Name, one data row
$wb.Worksheets.Item("_set_SCENARIOS").Range("A1").Copy($wb.Worksheets.Item("_set_EMISSIONS").Range("A1"))
$wsEmissions = $wb.Worksheets.Item("_set_EMISSIONS")
$wsEmissions.Range("A1").Value = "e_Name"
$wsEmissions.Range("A2").Value = "CO2 emissions"

# _set_LOADFACTORS: two-column header lf_Name / lf_Category, two data rows
$wb.Worksheets.Item("_set_FLOWS").Range("A1:B1").Copy($wb.Worksheets.Item("_set_LOADFACTORS").Range("A1:B1"))
$wsLoadFactors = $wb.Worksheets.Item("_set_LOADFACTORS")
$wsLoadFactors.Range("A1").Value = "lf_Name"
$wsLoadFactors.Range("B1").Value = "lf_Category"
$wsLoadFactors.Range("A2").Value = "lf max"
$wsLoadFactors.Range("B2").Value = "Max"
$wsLoadFactors.Range("A3").Value = "lf min"
$wsLoadFactors.Range("B3").Value = "Min"

# ---------------------------------------------------------------------------
# 5) Worksheet selections / active cells, per the target view state
#    (re-fetch every sheet by name since indices shifted after inserts)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("_set_TECH_SUP").Range("D15").Select()
$wb.Worksheets.Item("_set_TECH_DEM").Range("C8").Select()
$wb.Worksheets.Item("_set_FLOWS").Range("A2:B4").Select()
$wb.Worksheets.Item("_set_COSTS").Range("E21").Select()
$wb.Worksheets.Item("_set_EMISSIONS").Range("A2").Select()
$wb.Worksheets.Item("_set_LOADFACTORS").Range("G10").Select()
$wb.Worksheets.Item("_set_YEARS").Range("D1:D6").Select()

# ---------------------------------------------------------------------------
# 6) Active sheet / tab selection: _set_SCENARIOS tab is now the selected tab
#    and its active cell moves to J24.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("_set_SCENARIOS").Select()
$wb.Worksheets.Item("_set_SCENARIOS").Range("J24").Select()
